$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 61 - shifts existing rows 61-128 down to 62-129,
# carrying formatting (date style in column D) along with them.
$ws.Rows(61).Insert()

# Populate the new row 61 with the latest weekly price entry.
$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 44705
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 100112043
$ws.Range("G61").Value = "Pepino ensalada"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 150
$ws.Range("K61").Value = 17000
$ws.Range("L61").Value = 18000
$ws.Range("M61").Value = 17333
$ws.Range("N61").Value = "$/caja 60 unidades"
$ws.Range("O61").Value = "Región de Arica y Parinacota"
$ws.Range("P61").Value = 289
$ws.Range("Q61").Value = 60
$ws.Range("R61").Value = "Hortaliza"
